$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.901.21"
$ws.Range("E2").Value = "  +0.86%  "

# Row 3
$ws.Range("D3").Value = "2.734.55"
$ws.Range("E3").Value = "  +3.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'608.45"
$ws.Range("E5").Value = "  +1.87%  "

# Row 6
$ws.Range("D6").Value = "'169.66"
$ws.Range("E6").Value = "  +6.67%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").Value = "'0.549"
$ws.Range("E8").Value = "  +1.43%  "

# Row 9
$ws.Range("D9").Value = "2.735.07"
$ws.Range("E9").Value = "  +3.40%  "

# Row 10
$ws.Range("E10").Value = "  +4.08%  "

# Row 11
$ws.Range("D11").Value = "'0.368"
$ws.Range("E11").Value = "  +4.95%  "

# Row 12
$ws.Range("D12").Value = "'5.35"
$ws.Range("E12").Value = "  +1.08%  "

# Row 13
$ws.Range("E13").Value = "  -0.26%  "

# Row 14
$ws.Range("D14").Value = "'28.83"
$ws.Range("E14").Value = "  +3.15%  "

# Row 15
$ws.Range("D15").Value = "3.234.76"
$ws.Range("E15").Value = "  +3.39%  "

# Row 16
$ws.Range("D16").Value = "'0.0000191"
$ws.Range("E16").Value = "  +2.36%  "

# Row 17
$ws.Range("D17").Value = "68.974.65"
$ws.Range("E17").Value = "  +1.22%  "

# Row 18
$ws.Range("D18").Value = "2.688.42"
$ws.Range("E18").Value = "  +1.30%  "

# Row 19
$ws.Range("D19").Value = "'11.90"
$ws.Range("E19").Value = "  +4.80%  "

# Row 20
$ws.Range("D20").Value = "'377.18"
$ws.Range("E20").Value = "  +4.87%  "

# Row 21
$ws.Range("D21").Value = "'7.70"
$ws.Range("E21").Value = "  +3.99%  "

# Row 22
$ws.Range("D22").Value = "'4.53"
$ws.Range("E22").Value = "  +2.54%  "

# Row 23
$ws.Range("D23").Value = "'5.01"
$ws.Range("E23").Value = "  +5.47%  "

# Row 24
$ws.Range("D24").Value = "'2.10"
$ws.Range("E24").Value = "  +2.32%  "

# Row 25
$ws.Range("D25").Value = "'73.77"
$ws.Range("E25").Value = "  -0.99%  "

# Row 27
$ws.Range("D27").Value = "'10.16"
$ws.Range("E27").Value = "  +4.60%  "

# Row 28
$ws.Range("E28").Value = "  +3.47%  "

# Row 29
$ws.Range("D29").Value = "'0.0000106"
$ws.Range("E29").Value = "  +2.80%  "

# Row 30
$ws.Range("D30").Value = "'590.49"
$ws.Range("E30").Value = "  +5.30%  "

# Row 31
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
$ws.Range("D32").Value = "'8.35"
$ws.Range("E32").Value = "  +4.22%  "

# Row 33
$ws.Range("D33").Value = "'1.45"
$ws.Range("E33").Value = "  +4.25%  "

# Row 34
$ws.Range("D34").Value = "'1.99"
$ws.Range("E34").Value = "  +6.04%  "

# Row 35
$ws.Range("D35").Value = "'0.132"
$ws.Range("E35").Value = "  +3.85%  "

# Row 36
$ws.Range("E36").Value = "  -1.04%  "

# Row 37
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.14%  "

# Row 38
$ws.Range("D38").Value = "'163.08"
$ws.Range("E38").Value = "  +2.42%  "

# Row 39
$ws.Range("D39").Value = "'20.01"
$ws.Range("E39").Value = "  +1.71%  "

# Row 40
$ws.Range("D40").Value = "'0.382"
$ws.Range("E40").Value = "  +3.41%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.50"
$ws.Range("E41").Value = "  +3.27%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.91"
$ws.Range("E42").Value = "  +2.31%  "

# Row 43
$ws.Range("D43").Value = "'2.68"
$ws.Range("E43").Value = "  +2.53%  "

# Row 44
$ws.Range("D44").Value = "'17.99"
$ws.Range("E44").Value = "  +1.08%  "

# Row 45
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").Value = "0.0₆0313"
$ws.Range("E46").Value = "  -1.65%  "

# Row 47
$ws.Range("D47").Value = "'41.14"
$ws.Range("E47").Value = "  +1.67%  "

# Row 48
$ws.Range("D48").Value = "'0.607"
$ws.Range("E48").Value = "  +5.58%  "

# Row 49
$ws.Range("D49").Value = "'155.92"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50
$ws.Range("D50").Value = "'3.94"
$ws.Range("E50").Value = "  +3.68%  "

# Row 51
$ws.Range("E51").Value = "  +6.15%  "
